$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 129, pushing existing rows 129-133 down to 130-134.
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with the new weekly data point.
$ws.Cells.Item(129, 1).Value = 4
$ws.Cells.Item(129, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(129, 3).Value = "Los Lagos"
$ws.Cells.Item(129, 4).Value = 44509
$ws.Cells.Item(129, 5).Value = 10
$ws.Cells.Item(129, 6).Value = 100112039
$ws.Cells.Item(129, 7).Value = "Ciboulette"
$ws.Cells.Item(129, 8).Value = "Sin especificar"
$ws.Cells.Item(129, 9).Value = "Primera"
$ws.Cells.Item(129, 10).Value = 240
$ws.Cells.Item(129, 11).Value = 2500
$ws.Cells.Item(129, 12).Value = 2500
$ws.Cells.Item(129, 13).Value = 2500
$ws.Cells.Item(129, 14).Value = "$/docena de atados"
$ws.Cells.Item(129, 15).Value = "Región Metropolitana"
$ws.Cells.Item(129, 16).Value = 833
$ws.Cells.Item(129, 17).Value = 3
$ws.Cells.Item(129, 18).Value = "Hortaliza"
